$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 69.666664
$ws.Range("I4").Value = 69.666664
$ws.Range("K4").Value = 69.666664
$ws.Range("M4").Value = 44.333336

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 622.5325
$ws.Range("J17").Value = 564.2353000000001
$ws.Range("L17").Value = 1692.7059
$ws.Range("N17").Value = -2028.7059

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5283.7144
$ws.Range("J43").Value = 5998.3335
$ws.Range("L43").Value = 5998.3335
$ws.Range("N43").Value = -6136.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 13917.818
$ws.Range("I86").Value = 13099.5
$ws.Range("K86").Value = 13099.5
$ws.Range("M86").Value = -11976.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 13917.818
$ws.Range("I89").Value = 13099.5
$ws.Range("K89").Value = 65497.5
$ws.Range("M89").Value = -59881.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 146778.06
$ws.Range("I32").Value = 151793.5
$ws.Range("J32").Value = 34766.668
$ws.Range("K32").Value = 151793.5
$ws.Range("L32").Value = 34766.668
$ws.Range("M32").Value = -151506.5
$ws.Range("N32").Value = -35340.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7394.5347
$ws.Range("I74").Value = 4910.027
$ws.Range("K74").Value = 4910.027
$ws.Range("M74").Value = -4036.027

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7394.5347
$ws.Range("I77").Value = 4910.027
$ws.Range("K77").Value = 24550.135
$ws.Range("M77").Value = -20182.135

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 66948.44
$ws.Range("I20").Value = 95881.91
$ws.Range("J20").Value = 3294.8
$ws.Range("K20").Value = 95881.91
$ws.Range("L20").Value = 3294.8
$ws.Range("M20").Value = -95634.91
$ws.Range("N20").Value = -3788.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 28013.834
$ws.Range("J82").Value = 49000
$ws.Range("L82").Value = 49000
$ws.Range("N82").Value = -49766

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 28013.834
$ws.Range("J85").Value = 49000
$ws.Range("L85").Value = 49000
$ws.Range("N85").Value = -51652

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 11908.5
$ws.Range("I99").Value = 11908.5
$ws.Range("K99").Value = 11908.5
$ws.Range("M99").Value = -10410.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3674.0833
$ws.Range("I31").Value = 4158.9
$ws.Range("J31").Value = 1250
$ws.Range("K31").Value = 4158.9
$ws.Range("L31").Value = 1250
$ws.Range("M31").Value = -3863.9
$ws.Range("N31").Value = -1840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3674.0833
$ws.Range("I34").Value = 4158.9
$ws.Range("J34").Value = 1250
$ws.Range("K34").Value = 4158.9
$ws.Range("L34").Value = 1250
$ws.Range("M34").Value = -3956.9
$ws.Range("N34").Value = -1654

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 15000
$ws.Range("I69").Value = 15000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14251

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 15000
$ws.Range("I72").Value = 15000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -41256

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 74210.234
$ws.Range("I86").Value = 152548
$ws.Range("K86").Value = 152548
$ws.Range("M86").Value = -151425

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 74210.234
$ws.Range("I89").Value = 152548
$ws.Range("K89").Value = 762740
$ws.Range("M89").Value = -757124

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 15699.223
$ws.Range("J95").Value = 15699.223
$ws.Range("L95").Value = 15699.223
$ws.Range("N95").Value = -21191.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2584.4814
$ws.Range("I122").Value = 2455.739
$ws.Range("K122").Value = 7367.217000000001
$ws.Range("M122").Value = -4917.217000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1898.2941
$ws.Range("I33").Value = 79.14286
$ws.Range("J33").Value = 3171.7
$ws.Range("K33").Value = 474.85716
$ws.Range("L33").Value = 19030.2
$ws.Range("M33").Value = -191.85716
$ws.Range("N33").Value = -19596.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 29084.047
$ws.Range("I113").Value = 419.5
$ws.Range("J113").Value = 40549.867
$ws.Range("K113").Value = 1258.5
$ws.Range("L113").Value = 121649.601
$ws.Range("M113").Value = 911.5
$ws.Range("N113").Value = -125989.601

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 663.3
$ws.Range("I114").Value = 239.83333
$ws.Range("J114").Value = 1298.5
$ws.Range("K114").Value = 719.49999
$ws.Range("L114").Value = 3895.5
$ws.Range("M114").Value = 2534.50001
$ws.Range("N114").Value = -10403.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 368070.2
$ws.Range("I21").Value = 403377.7
$ws.Range("K21").Value = 403377.7
$ws.Range("M21").Value = -403204.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 368070.2
$ws.Range("I30").Value = 403377.7
$ws.Range("K30").Value = 403377.7
$ws.Range("M30").Value = -403272.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 26912.182
$ws.Range("J70").Value = 14522.667
$ws.Range("L70").Value = 14522.667
$ws.Range("N70").Value = -15062.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 26912.182
$ws.Range("J73").Value = 14522.667
$ws.Range("L73").Value = 14522.667
$ws.Range("N73").Value = -16394.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3911.4546
$ws.Range("I102").Value = 3802.6
$ws.Range("K102").Value = 3802.6
$ws.Range("M102").Value = -2180.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 4748
$ws.Range("I3").Value = 4748
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4748
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4636
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2765.8333
$ws.Range("I7").Value = 1650
$ws.Range("J7").Value = 4997.5
$ws.Range("K7").Value = 1650
$ws.Range("L7").Value = 4997.5
$ws.Range("M7").Value = -1538
$ws.Range("N7").Value = -5221.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H15").Value = 4748
$ws.Range("I15").Value = 4748
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 4748
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -4578
$ws.Range("N15").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2568.95
$ws.Range("I22").Value = 2235.9092
$ws.Range("J22").Value = 2976
$ws.Range("K22").Value = 2235.9092
$ws.Range("L22").Value = 2976
$ws.Range("M22").Value = -1940.9092
$ws.Range("N22").Value = -3566

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2568.95
$ws.Range("I27").Value = 2235.9092
$ws.Range("J27").Value = 2976
$ws.Range("K27").Value = 2235.9092
$ws.Range("L27").Value = 2976
$ws.Range("M27").Value = -2128.9092
$ws.Range("N27").Value = -3190

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12117.3
$ws.Range("I40").Value = 12117.3
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 12117.3
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -11981.3
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2765.8333
$ws.Range("I126").Value = 1650
$ws.Range("J126").Value = 4997.5
$ws.Range("K126").Value = 4950
$ws.Range("L126").Value = 14992.5
$ws.Range("M126").Value = -2480
$ws.Range("N126").Value = -19932.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7498.2383
$ws.Range("I136").Value = 3860.6667
$ws.Range("K136").Value = 11582.0001
$ws.Range("M136").Value = -9032.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 9666.125
$ws.Range("I2").Value = 2555
$ws.Range("J2").Value = 30999.5
$ws.Range("K2").Value = 2555
$ws.Range("L2").Value = 30999.5
$ws.Range("M2").Value = -2443
$ws.Range("N2").Value = -31223.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 49999.832
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 49999.832
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
